# Commit: "Changes done in Every Execution Classes and Write Code for Pass n Fail"
#
# Content-level changes in the sheet grid:
#   A3: "sri111@gmail.com"  ->  "gk030994@gmail.com"
#   B3: "bookstore"         ->  1234  (numeric)
#
# (Hyperlink relationship targets / rIds and the hyperlink style on A2/A3
# are left untouched - the underlying diff does not alter the .rels file.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "gk030994@gmail.com"
$ws.Range("B3").Value = 1234

# Move the active selection, as reflected in the saved view state.
$ws.Range("B11").Select()
